# Updates cryptocurrency price/volume figures to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.573.72"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.957.85"
$ws.Range("E3").Value = "  -0.22%  "
$ws.Range("E4").Value = "  +0.15%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "244.30"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.50%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.620"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.28%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "58.40"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +0.06%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.367"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.83%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "56.30"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.99%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0866"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +8.81%  "
$ws.Range("E12").Value = "  +1.24%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "22.01"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.25%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.830"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.33%  "
$ws.Range("D15").Value = "2.246.22"
$ws.Range("E15").Value = "  -0.07%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "13.68"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -2.04%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "5.24"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("D18").Value = "1.965.66"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "36.502.38"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "0.0₃0881"
$ws.Range("E20").Value = "  +3.38%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "70.07"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.56%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "230.37"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.55%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.08"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.59%  "
$ws.Range("E24").Value = "  +0.02%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.47"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("E26").Value = "  +1.64%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.42"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.58%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "162.38"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.95%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +10.25%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "19.64"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.86%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.118"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  +4.96%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.72"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.84%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.0645"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +5.19%  "
$ws.Range("E35").Value = "  -2.22%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "6.41"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +8.36%  "
$ws.Range("E37").Value = "  +0.05%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.77"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.96%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.19"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.44%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "3.03"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0998"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("E44").Value = "  -0.10%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "16.14"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.38%  "
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("D47").Value = "1.356.18"
$ws.Range("E47").Value = "  +1.15%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "88.56"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.38%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.22"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -4.23%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.83"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.18%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "46.18"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +4.73%  "
